# Apply the BOM_BME280_PCB TH update:
#  - rename sheet (date refresh 2024-12-24 -> 2025-02-24)
#  - refresh JLCPCB Price / Stock figures for rows 2-6 (C1, C2, R1/R2, SHT1/SHT2, U1)
#
# The JLCPCB Price / Stock columns (K/L) hold values that look numeric
# ("0.0019", "4896130", ...) but are stored as plain text in the workbook
# (shared-string cells, not numeric cells). A leading apostrophe forces
# Excel to keep the literal text instead of auto-coercing it to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the sheet/tab to reflect the new revision date.
$ws.Name = "BOM_BME280_PCB TH_2025-02-24"

# 2. Row 2 - CL05A105KA5NQNC / C1 (1uF cap): JLCPCB Stock changed.
$ws.Range("L2").Value = "'4896130"

# 3. Row 3 - 0402B104K500CT / C2 (100nF cap): JLCPCB Price & Stock changed.
$ws.Range("K3").Value = "'0.0019"
$ws.Range("L3").Value = "'282795"

# 4. Row 4 - 0402WGF4701TCE / R1,R2 (4.7k resistor): JLCPCB Price & Stock changed.
$ws.Range("K4").Value = "'0.0006"
$ws.Range("L4").Value = "'5158098"

# 5. Row 5 - HY-4P / SHT1,SHT2 (connector): JLCPCB Price changed; Stock column collapses.
$ws.Range("K5").Value = "'0.0687"
$ws.Range("L5").Value = "'4"

# 6. Row 6 - BME280 / U1: JLCPCB Price & Stock changed.
$ws.Range("K6").Value = "'3.471"
$ws.Range("L6").Value = "'4781"
